# Apply updated cryptocurrency price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "602.60", "69.119.09").
# Temporarily force the whole column to Text format so Excel does not silently
# reinterpret these strings as numbers (which would lose trailing zeros / precision
# or mangle multi-dot values), then restore the default "Normal" style afterwards
# so cell formatting matches the original workbook.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '69.119.09'
$ws.Range('D3').Value = '3.748.41'
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '602.60'
$ws.Range('D6').Value = '168.40'
$ws.Range('D7').Value = '3.747.68'
$ws.Range('D9').Value = '0.541'
$ws.Range('D13').Value = '38.19'
$ws.Range('D15').Value = '4.377.09'
$ws.Range('D16').Value = '3.749.83'
$ws.Range('D17').Value = '69.146.52'
$ws.Range('D20').Value = '17.21'
$ws.Range('D21').Value = '11.13'
$ws.Range('D22').Value = '494.15'
$ws.Range('D23').Value = '0.725'
$ws.Range('D24').Value = '0.0000152'
$ws.Range('D25').Value = '84.92'
$ws.Range('D27').Value = '12.34'
$ws.Range('D28').Value = '10.18'
$ws.Range('D31').Value = '2.48'
$ws.Range('D32').Value = '8.06'
$ws.Range('D33').Value = '31.67'
$ws.Range('D34').Value = '3.894.22'
$ws.Range('D36').Value = '3.683.24'
$ws.Range('D39').Value = '5.89'
$ws.Range('D43').Value = '433.72'
$ws.Range('D45').Value = '48.58'
$ws.Range('D46').Value = '8.53'
$ws.Range('D48').Value = '40.66'
$ws.Range('D49').Value = '141.27'
$ws.Range('D50').Value = '2.795.33'

$ws.Range('D2:D51').Style = 'Normal'

# The "Volume(1h)" column (E) holds percentage text with surrounding spaces
# (e.g. "  +1.17%  "), which Excel always keeps as plain text.
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('E7').Value = '  +0.92%  '
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  +3.35%  '
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('E21').Value = '  +20.46%  '
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('E24').Value = '  +6.31%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +1.59%  '
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('E42').Value = '  +4.68%  '
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('E44').Value = '  +1.49%  '
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('E51').Value = '  +0.69%  '
